$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) for rows 2-10: change the stored serial date
# value 45183 (2023-09-14) to 45184 (2023-09-15), wherever it currently has
# that value, leaving any other dates untouched.
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
